# Update vm_pu results table (rows 2-25) with new simulation values
# following the 380 kV case run (B column slack bus voltage 1.05 -> 1.02 p.u.)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.036368503616036
$ws.Cells.Item(2, 4).Value = 1.038649677590726
$ws.Cells.Item(2, 5).Value = 1.045562381905004
$ws.Cells.Item(2, 6).Value = 1.057154424233269
$ws.Cells.Item(2, 9).Value = 1.038918393331246
$ws.Cells.Item(2, 10).Value = 1.041477200329312
$ws.Cells.Item(2, 11).Value = 1.041437108288838
$ws.Cells.Item(2, 12).Value = 1.048330273227786
$ws.Cells.Item(2, 13).Value = 1.059890171113681
$ws.Cells.Item(2, 14).Value = 1.042956217412884

$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.037197220205669
$ws.Cells.Item(3, 4).Value = 1.039251831739385
$ws.Cells.Item(3, 5).Value = 1.046333088764091
$ws.Cells.Item(3, 6).Value = 1.058084535406938
$ws.Cells.Item(3, 9).Value = 1.039102914674604
$ws.Cells.Item(3, 10).Value = 1.041950460660396
$ws.Cells.Item(3, 11).Value = 1.041849868443189
$ws.Cells.Item(3, 12).Value = 1.048912551810137
$ws.Cells.Item(3, 13).Value = 1.060633766890124
$ws.Cells.Item(3, 14).Value = 1.043430149827922

$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.037734090527606
$ws.Cells.Item(4, 4).Value = 1.039641971630819
$ws.Cells.Item(4, 5).Value = 1.046832772919297
$ws.Cells.Item(4, 6).Value = 1.058687658469361
$ws.Cells.Item(4, 9).Value = 1.039221462109692
$ws.Cells.Item(4, 10).Value = 1.0422566566155
$ws.Cells.Item(4, 11).Value = 1.042116760911925
$ws.Cells.Item(4, 12).Value = 1.049289647828204
$ws.Cells.Item(4, 13).Value = 1.061115570652941
$ws.Cells.Item(4, 14).Value = 1.043736780616404

$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.03795994125452
$ws.Cells.Item(5, 4).Value = 1.039806105826423
$ws.Cells.Item(5, 5).Value = 1.047043074211383
$ws.Cells.Item(5, 6).Value = 1.05894151570007
$ws.Cells.Item(5, 9).Value = 1.039271095212772
$ws.Cells.Item(5, 10).Value = 1.042385371742701
$ws.Cells.Item(5, 11).Value = 1.042228915586639
$ws.Cells.Item(5, 12).Value = 1.049448254858361
$ws.Cells.Item(5, 13).Value = 1.061318274338243
$ws.Cells.Item(5, 14).Value = 1.043865678533848

$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.037997871359769
$ws.Cells.Item(6, 4).Value = 1.039833671630663
$ws.Cells.Item(6, 5).Value = 1.047078398397379
$ws.Cells.Item(6, 6).Value = 1.058984157235396
$ws.Cells.Item(6, 9).Value = 1.039279416832328
$ws.Cells.Item(6, 10).Value = 1.042406982976782
$ws.Cells.Item(6, 11).Value = 1.04224774406364
$ws.Cells.Item(6, 12).Value = 1.049474890083948
$ws.Cells.Item(6, 13).Value = 1.061352318125584
$ws.Cells.Item(6, 14).Value = 1.043887320458362

$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.037737107766983
$ws.Cells.Item(7, 4).Value = 1.039644164333043
$ws.Cells.Item(7, 5).Value = 1.046835582057973
$ws.Cells.Item(7, 6).Value = 1.058691049328791
$ws.Cells.Item(7, 9).Value = 1.039222126113029
$ws.Cells.Item(7, 10).Value = 1.042258376552814
$ws.Cells.Item(7, 11).Value = 1.042118259714049
$ws.Cells.Item(7, 12).Value = 1.049291766846794
$ws.Cells.Item(7, 13).Value = 1.061118278588876
$ws.Cells.Item(7, 14).Value = 1.043738502996225

$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.036648439741539
$ws.Cells.Item(8, 4).Value = 1.038853072618784
$ws.Cells.Item(8, 5).Value = 1.045822640801982
$ws.Cells.Item(8, 6).Value = 1.057468493712356
$ws.Cells.Item(8, 9).Value = 1.038980928806292
$ws.Cells.Item(8, 10).Value = 1.041637147399156
$ws.Cells.Item(8, 11).Value = 1.041576641328803
$ws.Cells.Item(8, 12).Value = 1.048526989298314
$ws.Cells.Item(8, 13).Value = 1.060141337712921
$ws.Cells.Item(8, 14).Value = 1.043116391625913

$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.034735005755329
$ws.Cells.Item(9, 4).Value = 1.037463018776066
$ws.Cells.Item(9, 5).Value = 1.044045334830309
$ws.Cells.Item(9, 6).Value = 1.055324073130622
$ws.Cells.Item(9, 9).Value = 1.038549427571286
$ws.Cells.Item(9, 10).Value = 1.040542248251617
$ws.Cells.Item(9, 11).Value = 1.040620831372945
$ws.Cells.Item(9, 12).Value = 1.047181889132225
$ws.Cells.Item(9, 13).Value = 1.0584248756132
$ws.Cells.Item(9, 14).Value = 1.042019937595999

$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.033462801171579
$ws.Cells.Item(10, 4).Value = 1.036539075599532
$ws.Cells.Item(10, 5).Value = 1.042865699232217
$ws.Cells.Item(10, 6).Value = 1.053901213811387
$ws.Cells.Item(10, 9).Value = 1.038257444417105
$ws.Cells.Item(10, 10).Value = 1.039812245792888
$ws.Cells.Item(10, 11).Value = 1.039982749225383
$ws.Cells.Item(10, 12).Value = 1.046286945494941
$ws.Cells.Item(10, 13).Value = 1.057284049484269
$ws.Cells.Item(10, 14).Value = 1.041288898450045

$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.032912753869701
$ws.Cells.Item(11, 4).Value = 1.036139674302218
$ws.Cells.Item(11, 5).Value = 1.042356168463004
$ws.Cells.Item(11, 6).Value = 1.053286725830584
$ws.Cells.Item(11, 9).Value = 1.038129997648608
$ws.Cells.Item(11, 10).Value = 1.039496145790144
$ws.Cells.Item(11, 11).Value = 1.039706260138893
$ws.Cells.Item(11, 12).Value = 1.045899867485384
$ws.Cells.Item(11, 13).Value = 1.056790905796714
$ws.Cells.Item(11, 14).Value = 1.040972349549041

$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.032708567513784
$ws.Cells.Item(12, 4).Value = 1.035991421636819
$ws.Cells.Item(12, 5).Value = 1.04216709708167
$ws.Cells.Item(12, 6).Value = 1.05305872286793
$ws.Cells.Item(12, 9).Value = 1.038082506275712
$ws.Cells.Item(12, 10).Value = 1.039378732813471
$ws.Cells.Item(12, 11).Value = 1.03960353165995
$ws.Cells.Item(12, 12).Value = 1.045756156803842
$ws.Cells.Item(12, 13).Value = 1.05660785838252
$ws.Cells.Item(12, 14).Value = 1.040854769832473

$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.032752360465343
$ws.Cells.Item(13, 4).Value = 1.036023217658315
$ws.Cells.Item(13, 5).Value = 1.042207644860457
$ws.Cells.Item(13, 6).Value = 1.053107619149164
$ws.Cells.Item(13, 9).Value = 1.038092700202137
$ws.Cells.Item(13, 10).Value = 1.039403918246431
$ws.Cells.Item(13, 11).Value = 1.039625568520358
$ws.Cells.Item(13, 12).Value = 1.045786980162986
$ws.Cells.Item(13, 13).Value = 1.056647116850145
$ws.Cells.Item(13, 14).Value = 1.040879991031637

$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.032895873187943
$ws.Cells.Item(14, 4).Value = 1.036127417584152
$ws.Cells.Item(14, 5).Value = 1.042340535855933
$ws.Cells.Item(14, 6).Value = 1.053267874020746
$ws.Cells.Item(14, 9).Value = 1.038126075092997
$ws.Cells.Item(14, 10).Value = 1.039486440376558
$ws.Cells.Item(14, 11).Value = 1.039697769137386
$ws.Cells.Item(14, 12).Value = 1.045887986933815
$ws.Cells.Item(14, 13).Value = 1.056775772417286
$ws.Cells.Item(14, 14).Value = 1.040962630352654

$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.032984312819246
$ws.Cells.Item(15, 4).Value = 1.036191632267059
$ws.Cells.Item(15, 5).Value = 1.04242243975324
$ws.Cells.Item(15, 6).Value = 1.053366644889583
$ws.Cells.Item(15, 9).Value = 1.038146618348273
$ws.Cells.Item(15, 10).Value = 1.039537285106774
$ws.Cells.Item(15, 11).Value = 1.039742250630605
$ws.Cells.Item(15, 12).Value = 1.045950229503609
$ws.Cells.Item(15, 13).Value = 1.056855058384127
$ws.Cells.Item(15, 14).Value = 1.041013547288219

$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.03349932353666
$ws.Cells.Item(16, 4).Value = 1.036565596861766
$ws.Cells.Item(16, 5).Value = 1.042899541799553
$ws.Cells.Item(16, 6).Value = 1.053942029639856
$ws.Cells.Item(16, 9).Value = 1.038265881289806
$ws.Cells.Item(16, 10).Value = 1.039833224329547
$ws.Cells.Item(16, 11).Value = 1.040001094884804
$ws.Cells.Item(16, 12).Value = 1.04631264394967
$ws.Cells.Item(16, 13).Value = 1.05731679568016
$ws.Cells.Item(16, 14).Value = 1.041309906778632

$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.033822598377286
$ws.Cells.Item(17, 4).Value = 1.036800356188078
$ws.Cells.Item(17, 5).Value = 1.043199153744549
$ws.Cells.Item(17, 6).Value = 1.05430338823337
$ws.Cells.Item(17, 9).Value = 1.038340420223264
$ws.Cells.Item(17, 10).Value = 1.040018858969821
$ws.Cells.Item(17, 11).Value = 1.040163409619819
$ws.Cells.Item(17, 12).Value = 1.046540095374968
$ws.Cells.Item(17, 13).Value = 1.057606657776552
$ws.Cells.Item(17, 14).Value = 1.041495805041388

$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.034011238677296
$ws.Cells.Item(18, 4).Value = 1.036937352069827
$ws.Cells.Item(18, 5).Value = 1.043374033696449
$ws.Cells.Item(18, 6).Value = 1.054514318715849
$ws.Cells.Item(18, 9).Value = 1.038383799432984
$ws.Cells.Item(18, 10).Value = 1.040127136032811
$ws.Cells.Item(18, 11).Value = 1.040258066132529
$ws.Cells.Item(18, 12).Value = 1.046672806154088
$ws.Cells.Item(18, 13).Value = 1.0577758105187
$ws.Cells.Item(18, 14).Value = 1.041604235870228

$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.034075573588311
$ws.Cells.Item(19, 4).Value = 1.036984075052103
$ws.Cells.Item(19, 5).Value = 1.043433683771549
$ws.Cells.Item(19, 6).Value = 1.054586266983099
$ws.Cells.Item(19, 9).Value = 1.038398573963775
$ws.Cells.Item(19, 10).Value = 1.040164055593151
$ws.Cells.Item(19, 11).Value = 1.040290338309885
$ws.Cells.Item(19, 12).Value = 1.046718064215891
$ws.Cells.Item(19, 13).Value = 1.057833500935053
$ws.Cells.Item(19, 14).Value = 1.041641207860579

$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.033787905805665
$ws.Cells.Item(20, 4).Value = 1.036775162025127
$ws.Cells.Item(20, 5).Value = 1.043166995650797
$ws.Cells.Item(20, 6).Value = 1.054264601713168
$ws.Cells.Item(20, 9).Value = 1.03833243304194
$ws.Cells.Item(20, 10).Value = 1.039998942170529
$ws.Cells.Item(20, 11).Value = 1.040145996733689
$ws.Cells.Item(20, 12).Value = 1.046515687613098
$ws.Cells.Item(20, 13).Value = 1.057575549917715
$ws.Cells.Item(20, 14).Value = 1.041475859957957

$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.032853608784784
$ws.Cells.Item(21, 4).Value = 1.036096730452235
$ws.Cells.Item(21, 5).Value = 1.042301397481694
$ws.Cells.Item(21, 6).Value = 1.053220676170215
$ws.Cells.Item(21, 9).Value = 1.038116251208975
$ws.Cells.Item(21, 10).Value = 1.039462139638311
$ws.Cells.Item(21, 11).Value = 1.039676508621712
$ws.Cells.Item(21, 12).Value = 1.045858241087183
$ws.Cells.Item(21, 13).Value = 1.056737883009538
$ws.Cells.Item(21, 14).Value = 1.040938295104571

$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.032266907265031
$ws.Cells.Item(22, 4).Value = 1.035670769006111
$ws.Cells.Item(22, 5).Value = 1.041758267731259
$ws.Cells.Item(22, 6).Value = 1.052565738547038
$ws.Cells.Item(22, 9).Value = 1.037979450212675
$ws.Cells.Item(22, 10).Value = 1.039124634899127
$ws.Cells.Item(22, 11).Value = 1.039381160713687
$ws.Cells.Item(22, 12).Value = 1.045445268606067
$ws.Cells.Item(22, 13).Value = 1.056211950604604
$ws.Cells.Item(22, 14).Value = 1.040600311069946

$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.032577859272718
$ws.Cells.Item(23, 4).Value = 1.035896522199743
$ws.Cells.Item(23, 5).Value = 1.042046085618701
$ws.Cells.Item(23, 6).Value = 1.052912798149668
$ws.Cells.Item(23, 9).Value = 1.038052054086094
$ws.Cells.Item(23, 10).Value = 1.039303551717642
$ws.Cells.Item(23, 11).Value = 1.039537745185189
$ws.Cells.Item(23, 12).Value = 1.045664155657868
$ws.Cells.Item(23, 13).Value = 1.056490686428154
$ws.Cells.Item(23, 14).Value = 1.040779481970864

$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.033803581639846
$ws.Cells.Item(24, 4).Value = 1.036786545983095
$ws.Cells.Item(24, 5).Value = 1.043181526134349
$ws.Cells.Item(24, 6).Value = 1.05428212718993
$ws.Cells.Item(24, 9).Value = 1.038336042408541
$ws.Cells.Item(24, 10).Value = 1.040007941716467
$ws.Cells.Item(24, 11).Value = 1.040153864926036
$ws.Cells.Item(24, 12).Value = 1.046526716299824
$ws.Cells.Item(24, 13).Value = 1.057589605969769
$ws.Cells.Item(24, 14).Value = 1.041484872284282

$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.035229079029728
$ws.Cells.Item(25, 4).Value = 1.037821902191335
$ws.Cells.Item(25, 5).Value = 1.04450389620507
$ws.Cells.Item(25, 6).Value = 1.055877274995038
$ws.Cells.Item(25, 9).Value = 1.038549427571286
$ws.Cells.Item(25, 10).Value = 1.04082532364481
$ws.Cells.Item(25, 11).Value = 1.040868090739263
$ws.Cells.Item(25, 12).Value = 1.047529320529386
$ws.Cells.Item(25, 13).Value = 1.058868015701874
$ws.Cells.Item(25, 14).Value = 1.04230341498872
